$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 135, pushing the existing rows 135-138 down to 137-140.
$ws.Rows.Item(135).Resize(2).Insert()

# New row 135: Ají, Inferno, Extra
$ws.Cells.Item(135, 1).Value = 4
$ws.Cells.Item(135, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(135, 3).Value = "Los Lagos"
$ws.Cells.Item(135, 4).Value = 44448
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 5).Value = 10
$ws.Cells.Item(135, 6).Value = 100112021
$ws.Cells.Item(135, 7).Value = "Ají"
$ws.Cells.Item(135, 8).Value = "Inferno"
$ws.Cells.Item(135, 9).Value = "Extra"
$ws.Cells.Item(135, 10).Value = 40
$ws.Cells.Item(135, 11).Value = 50000
$ws.Cells.Item(135, 12).Value = 50000
$ws.Cells.Item(135, 13).Value = 50000
$ws.Cells.Item(135, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(135, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(135, 16).Value = 4167
$ws.Cells.Item(135, 17).Value = 12
$ws.Cells.Item(135, 18).Value = "Hortaliza"

# New row 136: Ají, Inferno, Primera
$ws.Cells.Item(136, 1).Value = 4
$ws.Cells.Item(136, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(136, 3).Value = "Los Lagos"
$ws.Cells.Item(136, 4).Value = 44448
$ws.Cells.Item(136, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(136, 5).Value = 10
$ws.Cells.Item(136, 6).Value = 100112021
$ws.Cells.Item(136, 7).Value = "Ají"
$ws.Cells.Item(136, 8).Value = "Inferno"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 40
$ws.Cells.Item(136, 11).Value = 45000
$ws.Cells.Item(136, 12).Value = 45000
$ws.Cells.Item(136, 13).Value = 45000
$ws.Cells.Item(136, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(136, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(136, 16).Value = 3750
$ws.Cells.Item(136, 17).Value = 12
$ws.Cells.Item(136, 18).Value = "Hortaliza"
